$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Offer date
Replace-All "24-09-2025" "03-10-2025"

# Candidate name (appears multiple times throughout the letter)
Replace-All "Virat Kohli A" "Radhe Shyam"

# Address lines
Replace-All "Sudama Nagar," "Surya Nagar,"
Replace-All "Maha Laxmi Nagar, Maharashtra, 451111" "Maha Laxmi Nagar, Maharashtra, 452066"

# Contact details
Replace-All "9879809878" "8845679893"
Replace-All "Virat@gmail.com" "Radheshyam@gmail.com"

# Designation (appears multiple times)
Replace-All "Team Lead" "Solution Architect"

# Joining date
Replace-All "25-09-2025" "08-10-2025"

# CTC figures
Replace-All "16,00,000" "17,00,000"
Replace-All "Sixteen Lakh Rupees Only" "Seventeen Lakh Rupees Only"

# Compensation table figures
Replace-All "5,12,000" "5,44,000"
Replace-All "42,667" "45,333"
Replace-All "2,56,000" "2,72,000"
Replace-All "21,333" "22,667"
Replace-All "51,200" "54,400"
Replace-All "4,267" "4,533"
Replace-All "61,440" "65,280"
Replace-All "5,120" "5,440"
Replace-All "3,45,745" "3,80,598"
Replace-All "28,812" "31,716"
Replace-All "14,09,785" "14,99,678"
Replace-All "1,17,482" "1,24,973"
Replace-All "74,199" "78,930"
Replace-All "6,183" "6,578"
Replace-All "14,83,984" "15,78,608"
Replace-All "1,23,665" "1,31,551"
Replace-All "24,576" "26,112"
Replace-All "2,048" "2,176"
Replace-All "1,16,016" "1,21,392"
Replace-All "9,668" "10,116"
Replace-All "1,33,333" "1,41,667"
Replace-All "63,940" "67,780"
Replace-All "5,320" "5,640"
